$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) "Simulation Time" sheet: add new columns/rows of data + formulas
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Simulation Time")

# Row 1 headers (new shared strings SL 1..SL 10, PH 1..PH 10)
$ws2.Range("H1").Value = "SL 1"
$ws2.Range("I1").Value = "SL 2"
$ws2.Range("J1").Value = "SL 3"
$ws2.Range("K1").Value = "SL 4"
$ws2.Range("L1").Value = "SL 5"
$ws2.Range("M1").Value = "SL 6"
$ws2.Range("N1").Value = "SL 7"
$ws2.Range("O1").Value = "SL 8"
$ws2.Range("P1").Value = "SL 9"
$ws2.Range("Q1").Value = "SL 10"
$ws2.Range("S1").Value = "PH 1"
$ws2.Range("T1").Value = "PH 2"
$ws2.Range("U1").Value = "PH 3"
$ws2.Range("V1").Value = "PH 4"
$ws2.Range("W1").Value = "PH 5"
$ws2.Range("X1").Value = "PH 6"
$ws2.Range("Y1").Value = "PH 7"
$ws2.Range("Z1").Value = "PH 8"
$ws2.Range("AA1").Value = "PH 9"
$ws2.Range("AB1").Value = "PH 10"

# Column A row labels (new shared strings TTS, NHC, WH, MTG, NP)
$ws2.Range("A2").Value = "TTS"
$ws2.Range("A3").Value = "NHC"
$ws2.Range("A4").Value = "WH"
$ws2.Range("A5").Value = "MTG"
$ws2.Range("A6").Value = "NP"

# Row 3 raw trial data
$ws2.Range("H3").Value = 31717
$ws2.Range("I3").Value = 31548
$ws2.Range("J3").Value = 31389
$ws2.Range("K3").Value = 31221
$ws2.Range("L3").Value = 31595
$ws2.Range("M3").Value = 31288
$ws2.Range("N3").Value = 31576
$ws2.Range("O3").Value = 31467
$ws2.Range("P3").Value = 31612
$ws2.Range("Q3").Value = 31355
$ws2.Range("S3").Value = 749
$ws2.Range("T3").Value = 718
$ws2.Range("U3").Value = 702
$ws2.Range("V3").Value = 719
$ws2.Range("W3").Value = 718
$ws2.Range("X3").Value = 702
$ws2.Range("Y3").Value = 718
$ws2.Range("Z3").Value = 702
$ws2.Range("AA3").Value = 702
$ws2.Range("AB3").Value = 704

# Average formulas in columns B (Simulink) and C (Piha)
$ws2.Range("B2").Formula = "=AVERAGE(H2:Q2)"
$ws2.Range("C2").Formula = "=AVERAGE(S2:AB2)"
$ws2.Range("B3:B5").Formula = "=AVERAGE(H3:Q3)"
$ws2.Range("C3:C5").Formula = "=AVERAGE(S3:AB3)"
$ws2.Range("B6").Formula = "=AVERAGE(H6:Q6)"
$ws2.Range("C6").Formula = "=AVERAGE(S6:AB6)"

# ---------------------------------------------------------------------
# 2) Chart axis rescale: the "# Cells" value axis max 1000 -> 2000
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Scalability")
$chartObj = $ws1.ChartObjects(1)
$chart = $chartObj.Chart
$axis = $chart.Axes(1)
$axis.MaximumScale = 2000

# ---------------------------------------------------------------------
# 3) View state: active tab / selections on each sheet
# ---------------------------------------------------------------------
$ws1.Activate()
$ws1.Range("J50").Select()

$ws2.Activate()
$ws2.Range("D14").Select()
